$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Tests for form" section gains two new rows:
#   - a new "Avto" field test row (originally missing) right after the
#     "Positive" row, before the existing "Negative" rows.
#   - a new blank row (only a trailing "Pass") right after the existing
#     "Negative" rows, before the "User logout" header.
# Inserting shifts every row below down, which Excel also reflects in the
# existing merged ranges automatically.
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(12).Insert()

# New row 8: Avto field test (Name/Review/Rating style columns reused: here
# just C/D/E with "+" plus a trailing Pass result).
$ws.Cells.Item(8, 3).Value2 = "Avto"
$ws.Cells.Item(8, 4).Value2 = "+"
$ws.Cells.Item(8, 5).Value2 = "+"
$ws.Cells.Item(8, 6).Value2 = "Pass"

# New row 12: blank row, only the trailing Pass result is filled in.
$ws.Cells.Item(12, 6).Value2 = "Pass"

# Cursor moved to H9 in the saved file.
$ws.Range("H9").Select()
